# Applies the "Brightness + colour for frontend + backend" edit:
#  1. The "150mm=126mm 100mm=84mm" paragraph loses its second half and gains
#     a new "A=120mm B=170mm" measurement (as several separate runs), the
#     picture paragraph right after it is replaced by an empty paragraph that
#     just carries the (relocated) "_GoBack" bookmark, and a brand-new
#     paragraph "A                             B" is inserted after that.
#  2. The now-orphaned "_GoBack" bookmark that used to sit in the
#     "4x Långa Paneler Hållare" paragraph is removed (it moved up to the
#     new empty paragraph created in step 1).

$d = $word.ActiveDocument

# --- Step 1: locate the "150mm=126mm 100mm=84mm" paragraph and the picture
# paragraph that immediately follows it, then replace both of them (as one
# contiguous range) with three new paragraphs.
$measurePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "150mm=126mm*") {
        $measurePara = $candidate
        break
    }
}

$pictureIndex = $measurePara.Index + 1
$picturePara = $d.Paragraphs.Item($pictureIndex)

$targetRange = $d.Range($measurePara.Range.Start, $picturePara.Range.End)

$replacementXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>150mm=126mm</w:t></w:r><w:r><w:t xml:space="preserve">   </w:t></w:r><w:r><w:t>A=</w:t></w:r><w:r><w:t>120</w:t></w:r><w:r><w:t>mm</w:t></w:r><w:r><w:t xml:space="preserve"> B=170</w:t></w:r><w:r><w:t>mm</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:r><w:t>A                             B</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$targetRange.InsertXML($replacementXml)

# --- Step 2: strip the stale "_GoBack" bookmark out of the
# "4x Långa Paneler Hållare" paragraph (it used to mark the end of the
# document before the picture paragraph existed; that role now belongs to
# the empty paragraph inserted above).
$holderPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "4x L*nga Paneler H*llare*") {
        $holderPara = $candidate
        break
    }
}

$holderXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="30FB118E" w14:textId="77777777" w:rsidR="001100F6" w:rsidRPr="001100F6" w:rsidRDefault="001100F6" w:rsidP="001100F6"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>4x L&#229;nga Paneler H&#229;llare</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$holderPara.Range.InsertXML($holderXml)
